$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1608.4706
$ws.Range("I132").Value = 1396.5
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 4189.5
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -1659.5
$ws.Range("N132").Value = -20060

$ws.Range("H135").Value = 26694.21
$ws.Range("I135").Value = 408.30768
$ws.Range("J135").Value = 83647
$ws.Range("K135").Value = 3674.76912
$ws.Range("L135").Value = 752823
$ws.Range("M135").Value = -1139.76912
$ws.Range("N135").Value = -757893

$ws.Range("H137").Value = 2999
$ws.Range("I137").Value = 2998.5
$ws.Range("K137").Value = 8995.5
$ws.Range("M137").Value = -6445.5

$ws.Range("H138").Value = 1723.86
$ws.Range("I138").Value = 1187.175
$ws.Range("J138").Value = 2081.65
$ws.Range("K138").Value = 3561.525
$ws.Range("L138").Value = 6244.950000000001
$ws.Range("M138").Value = 1578.475
$ws.Range("N138").Value = -16524.95

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18523.703
$ws.Range("I32").Value = 4039.9421
$ws.Range("K32").Value = 4039.9421
$ws.Range("M32").Value = -3752.9421

$ws.Range("H61").Value = 1931.4762
$ws.Range("I61").Value = 1934.091
$ws.Range("J61").Value = 1928.6
$ws.Range("K61").Value = 1934.091
$ws.Range("L61").Value = 1928.6
$ws.Range("M61").Value = -1722.091
$ws.Range("N61").Value = -2352.6

$ws.Range("H74").Value = 1101.9
$ws.Range("I74").Value = 1091.5555
$ws.Range("K74").Value = 1091.5555
$ws.Range("M74").Value = -217.5554999999999

$ws.Range("H77").Value = 1101.9
$ws.Range("I77").Value = 1091.5555
$ws.Range("K77").Value = 5457.7775
$ws.Range("M77").Value = -1089.7775

$ws.Range("H132").Value = 1580.2222
$ws.Range("I132").Value = 1540.5938
$ws.Range("J132").Value = 1897.25
$ws.Range("K132").Value = 4621.7814
$ws.Range("L132").Value = 5691.75
$ws.Range("M132").Value = -2091.7814
$ws.Range("N132").Value = -10751.75

$ws.Range("H136").Value = 1931.4762
$ws.Range("I136").Value = 1934.091
$ws.Range("J136").Value = 1928.6
$ws.Range("K136").Value = 5802.272999999999
$ws.Range("L136").Value = 5785.799999999999
$ws.Range("M136").Value = -3252.272999999999
$ws.Range("N136").Value = -10885.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6495.2915
$ws.Range("I20").Value = 6061.1113
$ws.Range("J20").Value = 7797.8335
$ws.Range("K20").Value = 6061.1113
$ws.Range("L20").Value = 7797.8335
$ws.Range("M20").Value = -5814.1113
$ws.Range("N20").Value = -8291.833500000001

$ws.Range("H64").Value = 971.7692
$ws.Range("J64").Value = 1021.1818
$ws.Range("L64").Value = 1021.1818
$ws.Range("N64").Value = -1471.1818

$ws.Range("H67").Value = 971.7692
$ws.Range("J67").Value = 1021.1818
$ws.Range("L67").Value = 1021.1818
$ws.Range("N67").Value = -2581.1818

$ws.Range("H81").Value = 31666.334
$ws.Range("J81").Value = 31666.334
$ws.Range("L81").Value = 31666.334
$ws.Range("N81").Value = -33788.334

$ws.Range("H84").Value = 31666.334
$ws.Range("J84").Value = 31666.334
$ws.Range("L84").Value = 94999.00199999999
$ws.Range("N84").Value = -105607.002

$ws.Range("H134").Value = 2911.889
$ws.Range("I134").Value = 2565.2856
$ws.Range("J134").Value = 4125
$ws.Range("K134").Value = 7695.8568
$ws.Range("L134").Value = 12375
$ws.Range("M134").Value = -5160.8568
$ws.Range("N134").Value = -17445

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws.Range("H138").Value = 96932.63
$ws.Range("J138").Value = 100125.9
$ws.Range("L138").Value = 100125.9
$ws.Range("N138").Value = -110405.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1063.5385
$ws.Range("I58").Value = 1123.0555
$ws.Range("J58").Value = 929.625
$ws.Range("K58").Value = 1123.0555
$ws.Range("L58").Value = 929.625
$ws.Range("M58").Value = -920.0554999999999
$ws.Range("N58").Value = -1335.625

$ws.Range("H99").Value = 15651.417
$ws.Range("I99").Value = 17052.572
$ws.Range("J99").Value = 5843.3335
$ws.Range("K99").Value = 17052.572
$ws.Range("L99").Value = 5843.3335
$ws.Range("M99").Value = -15554.572
$ws.Range("N99").Value = -8839.333500000001

$ws.Range("H126").Value = 15651.417
$ws.Range("I126").Value = 17052.572
$ws.Range("J126").Value = 5843.3335
$ws.Range("K126").Value = 51157.716
$ws.Range("L126").Value = 17530.0005
$ws.Range("M126").Value = -48687.716
$ws.Range("N126").Value = -22470.0005

$ws.Range("H132").Value = 4881.4546
$ws.Range("I132").Value = 5744
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 17232
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -14702
$ws.Range("N132").Value = -8060

$ws.Range("H134").Value = 2992.8958
$ws.Range("I134").Value = 2905.0278
$ws.Range("J134").Value = 3256.5
$ws.Range("K134").Value = 8715.0834
$ws.Range("L134").Value = 9769.5
$ws.Range("M134").Value = -6180.0834
$ws.Range("N134").Value = -14839.5

$ws.Range("H136").Value = 1063.5385
$ws.Range("I136").Value = 1123.0555
$ws.Range("J136").Value = 929.625
$ws.Range("K136").Value = 3369.1665
$ws.Range("L136").Value = 2788.875
$ws.Range("M136").Value = -819.1664999999998
$ws.Range("N136").Value = -7888.875

$ws.Range("H141").Value = 286296
$ws.Range("J141").Value = 321512.5
$ws.Range("L141").Value = 321512.5
$ws.Range("N141").Value = -331872.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1506.6666
$ws.Range("I102").Value = 1630.3684
$ws.Range("K102").Value = 1630.3684
$ws.Range("M102").Value = -8.368400000000065

$ws.Range("H132").Value = 2255.2166
$ws.Range("I132").Value = 1658.8776
$ws.Range("J132").Value = 4911.636
$ws.Range("K132").Value = 4976.6328
$ws.Range("L132").Value = 14734.908
$ws.Range("M132").Value = -2446.6328
$ws.Range("N132").Value = -19794.908

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6539.1904
$ws.Range("I7").Value = 7939.1577
$ws.Range("J7").Value = 5382.696
$ws.Range("K7").Value = 7939.1577
$ws.Range("L7").Value = 5382.696
$ws.Range("M7").Value = -7827.1577
$ws.Range("N7").Value = -5606.696

$ws.Range("H61").Value = 80223.88
$ws.Range("I61").Value = 73351.28999999999
$ws.Range("J61").Value = 112296
$ws.Range("K61").Value = 73351.28999999999
$ws.Range("L61").Value = 112296
$ws.Range("M61").Value = -73149.28999999999
$ws.Range("N61").Value = -112700

$ws.Range("H113").Value = 80223.88
$ws.Range("I113").Value = 73351.28999999999
$ws.Range("J113").Value = 112296
$ws.Range("K113").Value = 73351.28999999999
$ws.Range("L113").Value = 112296
$ws.Range("M113").Value = -71181.28999999999
$ws.Range("N113").Value = -116636

$ws.Range("H126").Value = 6539.1904
$ws.Range("I126").Value = 7939.1577
$ws.Range("J126").Value = 5382.696
$ws.Range("K126").Value = 23817.4731
$ws.Range("L126").Value = 16148.088
$ws.Range("M126").Value = -21347.4731
$ws.Range("N126").Value = -21088.088

$ws.Range("H136").Value = 2253.4614
$ws.Range("I136").Value = 1982.5714
$ws.Range("J136").Value = 4623.75
$ws.Range("K136").Value = 5947.7142
$ws.Range("L136").Value = 13871.25
$ws.Range("M136").Value = -3397.7142
$ws.Range("N136").Value = -18971.25

$ws.Range("H137").Value = 67714.5
$ws.Range("I137").Value = 35000
$ws.Range("K137").Value = 35000
$ws.Range("M137").Value = -29900

$ws.Range("H138").Value = 65000
$ws.Range("J138").Value = 65000
$ws.Range("L138").Value = 65000
$ws.Range("N138").Value = -75280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1844.1177
$ws.Range("I126").Value = 1759.4375
$ws.Range("J126").Value = 3199
$ws.Range("K126").Value = 5278.3125
$ws.Range("L126").Value = 9597
$ws.Range("M126").Value = -2808.3125
$ws.Range("N126").Value = -14537

$ws.Range("H132").Value = 2803.8628
$ws.Range("I132").Value = 2885.6738
$ws.Range("K132").Value = 8657.0214
$ws.Range("M132").Value = -6127.0214

$ws.Range("H136").Value = 503.80554
$ws.Range("I136").Value = 503.80554
$ws.Range("K136").Value = 1511.41662
$ws.Range("M136").Value = 1038.58338

$ws.Range("H137").Value = 30047458
$ws.Range("J137").Value = 30047458
$ws.Range("L137").Value = 30047458
$ws.Range("N137").Value = -30057658
